$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.964.39'
$ws.Range('E2').Value = '  +4.08%  '
$ws.Range('D3').Value = '2.282.56'
$ws.Range('E3').Value = '  +4.54%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '251.74'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('E6').Value = '  +3.75%  '
$ws.Range('D7').Value = '72.10'
$ws.Range('E7').Value = '  +8.67%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.660'
$ws.Range('E9').Value = '  +15.37%  '
$ws.Range('D10').Value = '39.23'
$ws.Range('E10').Value = '  +7.71%  '
$ws.Range('D11').Value = '59.96'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').Value = '0.0975'
$ws.Range('E12').Value = '  +4.74%  '
$ws.Range('D13').Value = '7.45'
$ws.Range('E13').Value = '  +8.05%  '
$ws.Range('D14').Value = '0.106'
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('D15').Value = '2.625.04'
$ws.Range('E15').Value = '  +4.96%  '
$ws.Range('D16').Value = '14.92'
$ws.Range('E16').Value = '  +3.52%  '
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('D18').Value = '2.253.99'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('D19').Value = '42.912.87'
$ws.Range('E19').Value = '  +4.21%  '
$ws.Range('E20').Value = '  +7.74%  '
$ws.Range('E21').Value = '  +3.78%  '
$ws.Range('D22').Value = '73.43'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('D23').Value = '234.48'
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('D24').Value = '2.18'
$ws.Range('E24').Value = '  +6.68%  '
$ws.Range('D25').Value = '4.05'
$ws.Range('E25').Value = '  +7.61%  '
$ws.Range('D26').Value = '11.50'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').Value = '2.13'
$ws.Range('E30').Value = '  +5.54%  '
$ws.Range('D31').Value = '167.99'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '21.11'
$ws.Range('E32').Value = '  +4.26%  '
$ws.Range('D33').Value = '6.56'
$ws.Range('E33').Value = '  +14.24%  '
$ws.Range('D34').Value = '0.127'
$ws.Range('E34').Value = '  +3.82%  '
$ws.Range('D35').Value = '31.83'
$ws.Range('E35').Value = '  +30.16%  '
$ws.Range('E36').Value = '  +9.07%  '
$ws.Range('E37').Value = '  +3.96%  '
$ws.Range('D38').Value = '4.52'
$ws.Range('E38').Value = '  +13.76%  '
$ws.Range('D39').Value = '4.79'
$ws.Range('E39').Value = '  +5.24%  '
$ws.Range('E40').Value = '  +3.28%  '
$ws.Range('D41').Value = '13.47'
$ws.Range('E41').Value = '  +18.76%  '
$ws.Range('E42').Value = '  +5.29%  '
$ws.Range('D43').Value = '5.86'
$ws.Range('E43').Value = '  +6.66%  '
$ws.Range('D44').Value = '0.212'
$ws.Range('E44').Value = '  +11.50%  '
$ws.Range('D45').Value = '9.21'
$ws.Range('E45').Value = '  +7.80%  '
$ws.Range('D46').Value = '62.58'
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('D47').Value = '5.00'
$ws.Range('E47').Value = '  -5.64%  '
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('E49').Value = '  +3.76%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  +4.25%  '
